$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 (Thang 7) updates
$ws.Range("D7").Value = "2024-07-17T12:15:00.000Z"
$ws.Range("T7").Value = 38900000
$ws.Range("W7").Value = 185471000
$ws.Range("AA7").Value = 141397000
$ws.Range("AE7").Value = 326868000
$ws.Range("AH7").Value = 275868000
$ws.Range("AK7").Value = 48
$ws.Range("AN7").Value = 51000000
$ws.Range("AQ7").Value = 314768000

# Row 8 (Thang 6) updates
$ws.Range("D8").Value = "2024-07-17T12:15:00.000Z"

# Row 9 (Thang 5) updates
$ws.Range("D9").Value = "2024-07-17T12:15:00.000Z"

# Row 10 (Thang 4) updates
$ws.Range("D10").Value = "2024-07-17T12:15:00.000Z"

# Row 11 (Thang 3) updates
$ws.Range("D11").Value = "2024-07-17T12:15:00.000Z"

# Row 12 (Thang 2) updates
$ws.Range("D12").Value = "2024-07-17T12:15:00.000Z"
$ws.Range("T12").Value = 49498000
$ws.Range("AK12").Value = 42
$ws.Range("AQ12").Value = 399293000
